$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 45 and 46 swap place in the ranking (EnergySwap <-> PancakeSwap),
# so their Coin name / Link / Price / Volume values are exchanged.
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.755"
$ws.Range("E45").Value = "  +0.86%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.01"
$ws.Range("E46").Value = "  -4.55%  "

# Price / Volume(1h) updates for the remaining rows
$ws.Range("D2").Value = "22.020.29"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.552.25"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.63"
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3924"
$ws.Range("E7").Value = "  +3.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3207"
$ws.Range("E8").Value = "  -1.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.32"
$ws.Range("E9").Value = "  +2.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07183"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.070"
$ws.Range("E11").Value = "  -4.96%  "
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.638"
$ws.Range("E13").Value = "  -2.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.58"
$ws.Range("E14").Value = "  -7.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.641"
$ws.Range("E15").Value = "  -1.60%  "
$ws.Range("D16").Value = "1.554.23"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001093"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06558"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("E19").Value = "  -2.62%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.181"
$ws.Range("E21").Value = "  -3.06%  "
$ws.Range("E22").Value = "  -3.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.13"
$ws.Range("E23").Value = "  -4.34%  "
$ws.Range("D24").Value = "22.043.50"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.368"
$ws.Range("E25").Value = "  +3.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.378"
$ws.Range("E26").Value = "  -4.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "147.94"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.46"
$ws.Range("E28").Value = "  -3.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.858"
$ws.Range("D30").Value = "1.729.26"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.08"
$ws.Range("E31").Value = "  -2.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9727"
$ws.Range("E32").Value = "  -9.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.785"
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("E34").Value = "  +1.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.099"
$ws.Range("E35").Value = "  -1.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.606"
$ws.Range("E36").Value = "  -13.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02246"
$ws.Range("E37").Value = "  -2.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.088"
$ws.Range("E38").Value = "  -2.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05981"
$ws.Range("E39").Value = "  -3.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.204"
$ws.Range("E40").Value = "  -2.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2049"
$ws.Range("E41").Value = "  -4.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.68"
$ws.Range("E43").Value = "  -2.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5787"
$ws.Range("E44").Value = "  -3.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5536"
$ws.Range("E47").Value = "  -4.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "116.98"
$ws.Range("E48").Value = "  -3.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.870"
$ws.Range("E49").Value = "  -4.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.131"
$ws.Range("E50").Value = "  -3.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06803"
$ws.Range("E51").Value = "  -2.93%  "
